$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Amazon account credentials used in rows 2 and 3
$ws.Range("A2").Value = "testdemo372@gmail.com"
$ws.Range("B2").Value = "india123"
$ws.Range("A3").Value = "testdemo372@gmail.com"
$ws.Range("B3").Value = "india123"

# Re-fit column A now that its content changed
$ws.Columns.Item(1).AutoFit()

# Move the active selection as it was left when the sheet was saved
$ws.Range("C7").Select()
